$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all the new / changed text values first, in row order, so the shared
# string table is built in the same order as the target file.
$ws.Range("A25").Value = "Background data @ -10m"
$ws.Range("A26").Value = "%diff in CO2 ppm during injection"
$ws.Range("A27").Value = "Ave CO2 ppm (air)"
$ws.Range("A28").Value = "Ave CO2 ppm (water)"
$ws.Range("A29").Value = "Ave Temp c (air)"
$ws.Range("B29").Value = 8.194
$ws.Range("A30").Value = "Ave pressure kpa (air)"
$ws.Range("B30").Value = 63.1
$ws.Range("A31").Value = "Ave Temp c (water)"
$ws.Range("B31").Value = 7.005
$ws.Range("A32").Value = "Ave flux um/m2"
$ws.Range("A33").Value = "Ave k m/d"
$ws.Range("A34").Value = "Ave k600 m/d"
$ws.Range("A35").Value = "ER"
$ws.Range("A36").Value = "GPP"

# Now apply bold formatting. A25 gets bolded (and font reset to a
# scheme-less Calibri) first, producing the "no-scheme" bold font that ends
# up as the most-used bold style. A27 is bolded next using the default
# scheme font, producing the second, less-used bold style.
$ws.Range("A25").Font.Bold = $true
$ws.Range("A25").Font.Name = "Calibri"

$ws.Range("A27").Font.Bold = $true

$ws.Range("A26").Font.Bold = $true
$ws.Range("A26").Font.Name = "Calibri"

$ws.Range("A28").Font.Bold = $true
$ws.Range("A28").Font.Name = "Calibri"

$ws.Range("A29").Font.Bold = $true
$ws.Range("A29").Font.Name = "Calibri"

$ws.Range("A30").Font.Bold = $true
$ws.Range("A30").Font.Name = "Calibri"

$ws.Range("A31").Font.Bold = $true
$ws.Range("A31").Font.Name = "Calibri"

$ws.Range("A32").Font.Bold = $true
$ws.Range("A32").Font.Name = "Calibri"

$ws.Range("A33").Font.Bold = $true
$ws.Range("A33").Font.Name = "Calibri"

$ws.Range("A34").Font.Bold = $true
$ws.Range("A34").Font.Name = "Calibri"

$ws.Range("A35").Font.Bold = $true
$ws.Range("A35").Font.Name = "Calibri"

$ws.Range("A36").Font.Bold = $true
$ws.Range("A36").Font.Name = "Calibri"

# Column A width
$ws.Columns.Item(1).ColumnWidth = 28.83203125

# Selection
$ws.Range("B31").Select()
